$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain plain text so values like "314.33"
# or "28.211.76" are not auto-coerced into numeric cells by Excel.
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "28.211.76"
$ws.Range("E2").Value = "  -0.68%  "
Set-TextCell $ws.Range("D3") "1.801.20"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextCell $ws.Range("D5") "314.33"
$ws.Range("E5").Value = "  -0.35%  "
Set-TextCell $ws.Range("D6") "1.002"
$ws.Range("E6").Value = "  +0.05%  "
Set-TextCell $ws.Range("D7") "0.5262"
$ws.Range("E7").Value = "  +3.37%  "
Set-TextCell $ws.Range("D8") "0.3820"
$ws.Range("E8").Value = "  -2.44%  "
Set-TextCell $ws.Range("D9") "0.07972"
$ws.Range("E9").Value = "  +3.04%  "
Set-TextCell $ws.Range("D10") "41.39"
$ws.Range("E10").Value = "  -1.00%  "
Set-TextCell $ws.Range("D11") "1.098"
$ws.Range("E11").Value = "  -0.82%  "
Set-TextCell $ws.Range("D12") "6.312"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  +0.12%  "
Set-TextCell $ws.Range("D14") "20.64"
$ws.Range("E14").Value = "  -1.57%  "
Set-TextCell $ws.Range("D15") "1.808.07"
$ws.Range("E15").Value = "  -0.95%  "
Set-TextCell $ws.Range("D16") "7.309"
$ws.Range("E16").Value = "  -2.32%  "
Set-TextCell $ws.Range("D17") "92.85"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("E18").Value = "  -4.88%  "
Set-TextCell $ws.Range("D19") "0.06605"
$ws.Range("E19").Value = "  -0.27%  "
Set-TextCell $ws.Range("D20") "1.002"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  -2.11%  "
Set-TextCell $ws.Range("D22") "5.964"
$ws.Range("E22").Value = "  -2.01%  "
Set-TextCell $ws.Range("D23") "28.249.97"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("E24").Value = "  -0.41%  "
Set-TextCell $ws.Range("D25") "2.236"
$ws.Range("E25").Value = "  -0.68%  "
Set-TextCell $ws.Range("D26") "159.39"
$ws.Range("E26").Value = "  +2.48%  "
Set-TextCell $ws.Range("D27") "20.50"
$ws.Range("E27").Value = "  -2.94%  "
Set-TextCell $ws.Range("D28") "2.006.84"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("E29").Value = "  -0.14%  "
Set-TextCell $ws.Range("D30") "123.12"
$ws.Range("E30").Value = "  -1.48%  "
Set-TextCell $ws.Range("D31") "0.1094"
$ws.Range("E31").Value = "  -0.47%  "
Set-TextCell $ws.Range("D32") "1.059"
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("E33").Value = "  +0.21%  "
Set-TextCell $ws.Range("D34") "5.531"
$ws.Range("E34").Value = "  -2.23%  "
Set-TextCell $ws.Range("D35") "0.07292"
$ws.Range("E35").Value = "  +3.47%  "
Set-TextCell $ws.Range("D36") "12.22"
$ws.Range("E36").Value = "  +9.44%  "

# Rows 37/38: swap Algorand and FraxShare, with updated prices
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws.Range("D37") "8.856"
$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell $ws.Range("D38") "0.2160"
$ws.Range("E38").Value = "  -2.31%  "

Set-TextCell $ws.Range("D39") "0.02305"
$ws.Range("E39").Value = "  -0.81%  "
Set-TextCell $ws.Range("D40") "5.067"
$ws.Range("E40").Value = "  -2.24%  "
Set-TextCell $ws.Range("D41") "0.6196"
$ws.Range("E41").Value = "  -1.04%  "
Set-TextCell $ws.Range("D42") "1.164"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("E43").Value = "  -1.38%  "
Set-TextCell $ws.Range("D44") "13.23"
$ws.Range("E44").Value = "  -1.43%  "
Set-TextCell $ws.Range("D45") "0.5989"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("E46").Value = "  +0.84%  "
Set-TextCell $ws.Range("D47") "126.80"
$ws.Range("E47").Value = "  +2.11%  "
Set-TextCell $ws.Range("D48") "1.206"
$ws.Range("E48").Value = "  +1.04%  "
Set-TextCell $ws.Range("D49") "1.921"
$ws.Range("E49").Value = "  -2.92%  "
Set-TextCell $ws.Range("D50") "0.06826"
$ws.Range("E50").Value = "  -1.01%  "
Set-TextCell $ws.Range("D51") "72.98"
$ws.Range("E51").Value = "  -1.47%  "
